{"js": "const pairs = [\n  [\"2025-05-10 Saturday\", \"2025-05-11 Sunday\"],\n  [\"200\u00f78=25, 0\", \"720\u00f76=120, 0\"],\n  [\"173\u00f78=21, 5\", \"627\u00f75=125, 2\"],\n  [\"775\u00f75=155, 0\", \"591\u00f74=147, 3\"],\n  [\"997\u00f75=199, 2\", \"897\u00f75=179, 2\"],\n  [\"509\u00f74=127, 1\", \"850\u00f75=170, 0\"],\n  [\"186\u00f73=62, 0\", \"720\u00f74=180, 0\"],\n  [\"417\u00f74=104, 1\", \"152\u00f79=16, 8\"],\n  [\"175\u00f77=25, 0\", \"606\u00f74=151, 2\"],\n  [\"239\u00f72=119, 1\", \"534\u00f79=59, 3\"],\n  [\"258\u00f77=36, 6\", \"658\u00f74=164, 2\"],\n  [\"728\u00f75=145, 3\", \"897\u00f79=99, 6\"],\n  [\"711\u00f73=237, 0\", \"460\u00f79=51, 1\"],\n  [\"491\u00f75=98, 1\", \"737\u00f74=184, 1\"],\n  [\"732\u00f73=244, 0\", \"617\u00f74=154, 1\"],\n  [\"116\u00f77=16, 4\", \"186\u00f72=93, 0\"],\n  [\"393\u00f77=56, 1\", \"566\u00f79=62, 8\"],\n  [\"261\u00f75=52, 1\", \"858\u00f77=122, 4\"],\n  [\"389\u00f79=43, 2\", \"465\u00f77=66, 3\"],\n  [\"288\u00f73=96, 0\", \"276\u00f75=55, 1\"],\n  [\"905\u00f76=150, 5\", \"279\u00f72=139, 1\"],\n  [\"154\u00f75=30, 4\", \"854\u00f77=122, 0\"],\n  [\"314\u00f76=52, 2\", \"719\u00f76=119, 5\"],\n  [\"845\u00f74=211, 1\", \"359\u00f75=71, 4\"],\n  [\"208\u00f78=26, 0\", \"863\u00f79=95, 8\"],\n  [\"638\u00f75=127, 3\", \"430\u00f76=71, 4\"],\n];\n\nconst body = context.document.body;\nlet totalMatches = 0;\n\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n    totalMatches++;\n  }\n  await context.sync();\n}\n\nreturn totalMatches;\n", "ps1": "$pairs = @(\n    @(\"2025-05-10 Saturday\", \"2025-05-11 Sunday\"),\n    @(\"200\u00f78=25, 0\", \"720\u00f76=120, 0\"),\n    @(\"173\u00f78=21, 5\", \"627\u00f75=125, 2\"),\n    @(\"775\u00f75=155, 0\", \"591\u00f74=147, 3\"),\n    @(\"997\u00f75=199, 2\", \"897\u00f75=179, 2\"),\n    @(\"509\u00f74=127, 1\", \"850\u00f75=170, 0\"),\n    @(\"186\u00f73=62, 0\", \"720\u00f74=180, 0\"),\n    @(\"417\u00f74=104, 1\", \"152\u00f79=16, 8\"),\n    @(\"175\u00f77=25, 0\", \"606\u00f74=151, 2\"),\n    @(\"239\u00f72=119, 1\", \"534\u00f79=59, 3\"),\n    @(\"258\u00f77=36, 6\", \"658\u00f74=164, 2\"),\n    @(\"728\u00f75=145, 3\", \"897\u00f79=99, 6\"),\n    @(\"711\u00f73=237, 0\", \"460\u00f79=51, 1\"),\n    @(\"491\u00f75=98, 1\", \"737\u00f74=184, 1\"),\n    @(\"732\u00f73=244, 0\", \"617\u00f74=154, 1\"),\n    @(\"116\u00f77=16, 4\", \"186\u00f72=93, 0\"),\n    @(\"393\u00f77=56, 1\", \"566\u00f79=62, 8\"),\n    @(\"261\u00f75=52, 1\", \"858\u00f77=122, 4\"),\n    @(\"389\u00f79=43, 2\", \"465\u00f77=66, 3\"),\n    @(\"288\u00f73=96, 0\", \"276\u00f75=55, 1\"),\n    @(\"905\u00f76=150, 5\", \"279\u00f72=139, 1\"),\n    @(\"154\u00f75=30, 4\", \"854\u00f77=122, 0\"),\n    @(\"314\u00f76=52, 2\", \"719\u00f76=119, 5\"),\n    @(\"845\u00f74=211, 1\", \"359\u00f75=71, 4\"),\n    @(\"208\u00f78=26, 0\", \"863\u00f79=95, 8\"),\n    @(\"638\u00f75=127, 3\", \"430\u00f76=71, 4\"),\n)\n\n$d = $word.ActiveDocument\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
